$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='62.395.11'; E='  +2.85%  '},
    @{Row=3; D='3.430.57'; E='  +1.87%  '},
    @{Row=4; E='  +0.10%  '},
    @{Row=5; D='406.69'; E='  +0.64%  '},
    @{Row=6; D='132.83'; E='  +4.73%  '},
    @{Row=7; D='0.595'; E='  -1.31%  '},
    @{Row=8; E='  -0.07%  '},
    @{Row=9; D='0.689'; E='  +3.38%  '},
    @{Row=10; D='0.135'; E='  +9.02%  '},
    @{Row=11; D='41.92'; E='  +0.38%  '},
    @{Row=12; E='  -0.14%  '},
    @{Row=13; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='8.40'; E='  -0.76%  '},
    @{Row=14; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='19.84'; E='  +1.62%  '},
    @{Row=15; D='3.438.11'; E='  +2.37%  '},
    @{Row=16; D='11.74'; E='  +3.58%  '},
    @{Row=17; D='62.279.17'; E='  +2.61%  '},
    @{Row=18; E='  +0.35%  '},
    @{Row=19; D='0.0000149'; E='  +11.91%  '},
    @{Row=20; E='  -2.28%  '},
    @{Row=21; D='84.07'; E='  +2.53%  '},
    @{Row=22; D='312.00'; E='  +2.60%  '},
    @{Row=23; D='12.77'; E='  -2.11%  '},
    @{Row=24; D='3.18'; E='  +1.12%  '},
    @{Row=25; D='4.74'; E='  +0.24%  '},
    @{Row=26; E='  +1.04%  '},
    @{Row=27; D='8.14'; E='  -4.20%  '},
    @{Row=28; E='  +5.26%  '},
    @{Row=29; D='2.74'; E='  +5.81%  '},
    @{Row=30; D='0.173'; E='  +0.11%  '},
    @{Row=31; D='43.80'; E='  +3.20%  '},
    @{Row=32; E='  -0.45%  '},
    @{Row=33; E='  -2.97%  '},
    @{Row=34; E='  -0.03%  '},
    @{Row=35; E='  +0.93%  '},
    @{Row=36; D='51.63'; E='  -0.91%  '},
    @{Row=37; D='0.999'; E='  +0.00%  '},
    @{Row=38; D='3.00'; E='  +1.79%  '},
    @{Row=39; E='  -2.08%  '},
    @{Row=40; D='0.318'; E='  +12.90%  '},
    @{Row=41; D='142.07'; E='  +4.63%  '},
    @{Row=42; E='  +0.61%  '},
    @{Row=43; E='  -3.83%  '},
    @{Row=44; D='3.92'; E='  +0.08%  '},
    @{Row=45; D='16.77'; E='  -0.42%  '},
    @{Row=46; E='  +0.43%  '},
    @{Row=47; D='21.28'; E='  -2.31%  '},
    @{Row=48; D='2.103.13'; E='  -1.42%  '},
    @{Row=49; D='2.32'; E='  -1.10%  '},
    @{Row=50; D='1.95'; E='  +3.32%  '},
    @{Row=51; E='  +19.62%  '}
)

foreach ($item in $updates) {
    $r = $item.Row
    if ($item.ContainsKey('B')) { $ws.Range("B$r").Value = $item.B }
    if ($item.ContainsKey('C')) { $ws.Range("C$r").Value = $item.C }
    if ($item.ContainsKey('D')) {
        $cell = $ws.Range("D$r")
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = $origStyle
    }
    if ($item.ContainsKey('E')) { $ws.Range("E$r").Value = $item.E }
}